# Applies translation-related changes to the follow_map_time workbook:
#  - survey sheet header F1: "display.text" -> "display.prompt.text"
#  - settings sheet header C1: "display.title" -> "display.title.text"
#  - updates active cell selections on survey/settings sheets
#  - makes the "survey" sheet the active tab (instead of "properties")

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")
$properties = $wb.Worksheets.Item("properties")

# Update the translated header labels (order matters for shared-string
# table layout: "display.title.text" is introduced before
# "display.prompt.text").
$settings.Range("C1").Value = "display.title.text"
$survey.Range("F1").Value = "display.prompt.text"

# Update selections left behind on each sheet.
$survey.Range("F2").Select()
$settings.Range("C2").Select()
$properties.Range("E5").Select()

# Make "survey" the active sheet/tab (was "properties").
$survey.Activate()
